# Scheduled-runner update: refresh cached market-board price/profit figures
# (currentAveragePrice*, Leve*Price*, Leve*Profit*) across the ALC, ARM, CRP,
# CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1753.7142
$ws.Range("I96").Value = 1145.2
$ws.Range("J96").Value = 3275
$ws.Range("K96").Value = 3435.6
$ws.Range("L96").Value = 9825
$ws.Range("M96").Value = -2062.6
$ws.Range("N96").Value = -12571

$ws.Range("H129").Value = 1276.625
$ws.Range("J129").Value = 1301.6774
$ws.Range("L129").Value = 3905.0322
$ws.Range("N129").Value = -13905.0322

$ws.Range("H132").Value = 3724.6365
$ws.Range("I132").Value = 4114.0586
$ws.Range("K132").Value = 12342.1758
$ws.Range("M132").Value = -9812.1758

$ws.Range("H137").Value = 101850.4
$ws.Range("I137").Value = 1750.5
$ws.Range("K137").Value = 5251.5
$ws.Range("M137").Value = -2701.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 20970.6
$ws.Range("J55").Value = 23713.25
$ws.Range("L55").Value = 23713.25
$ws.Range("N55").Value = -24343.25

$ws.Range("H63").Value = 3127140
$ws.Range("I63").Value = 2425
$ws.Range("K63").Value = 2425
$ws.Range("M63").Value = -1739

$ws.Range("H66").Value = 3127140
$ws.Range("I66").Value = 2425
$ws.Range("K66").Value = 12125
$ws.Range("M66").Value = -8693

$ws.Range("H80").Value = 47287
$ws.Range("J80").Value = 47287
$ws.Range("L80").Value = 47287
$ws.Range("N80").Value = -49283

$ws.Range("H83").Value = 47287
$ws.Range("J83").Value = 47287
$ws.Range("L83").Value = 141861
$ws.Range("N83").Value = -151845

$ws.Range("H88").Value = 86639.836
$ws.Range("I88").Value = 1640.4
$ws.Range("J88").Value = 147353.72
$ws.Range("K88").Value = 1640.4
$ws.Range("L88").Value = 147353.72
$ws.Range("M88").Value = -1234.4
$ws.Range("N88").Value = -148165.72

$ws.Range("H91").Value = 86639.836
$ws.Range("I91").Value = 1640.4
$ws.Range("J91").Value = 147353.72
$ws.Range("K91").Value = 1640.4
$ws.Range("L91").Value = 147353.72
$ws.Range("M91").Value = -236.4000000000001
$ws.Range("N91").Value = -150161.72

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1180.8334
$ws.Range("I16").Value = 1495
$ws.Range("J16").Value = 1023.75
$ws.Range("K16").Value = 1495
$ws.Range("L16").Value = 1023.75
$ws.Range("M16").Value = -1208
$ws.Range("N16").Value = -1597.75

$ws.Range("H58").Value = 16222.546
$ws.Range("I58").Value = 1117.32
$ws.Range("J58").Value = 63426.375
$ws.Range("K58").Value = 1117.32
$ws.Range("L58").Value = 63426.375
$ws.Range("M58").Value = -914.3199999999999
$ws.Range("N58").Value = -63832.375

$ws.Range("H107").Value = 1332.1428
$ws.Range("I107").Value = 659
$ws.Range("K107").Value = 659
$ws.Range("M107").Value = 1261

$ws.Range("H113").Value = 1180.8334
$ws.Range("I113").Value = 1495
$ws.Range("J113").Value = 1023.75
$ws.Range("K113").Value = 1495
$ws.Range("L113").Value = 1023.75
$ws.Range("M113").Value = 675
$ws.Range("N113").Value = -5363.75

$ws.Range("H132").Value = 31708.945
$ws.Range("I132").Value = 41488.23
$ws.Range("J132").Value = 6282.8
$ws.Range("K132").Value = 124464.69
$ws.Range("L132").Value = 18848.4
$ws.Range("M132").Value = -121934.69
$ws.Range("N132").Value = -23908.4

$ws.Range("H134").Value = 1123.129
$ws.Range("I134").Value = 927.5333
$ws.Range("J134").Value = 1306.5
$ws.Range("K134").Value = 2782.5999
$ws.Range("L134").Value = 3919.5
$ws.Range("M134").Value = -247.5999000000002
$ws.Range("N134").Value = -8989.5

$ws.Range("H136").Value = 16222.546
$ws.Range("I136").Value = 1117.32
$ws.Range("J136").Value = 63426.375
$ws.Range("K136").Value = 3351.96
$ws.Range("L136").Value = 190279.125
$ws.Range("M136").Value = -801.96
$ws.Range("N136").Value = -195379.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 742.44
$ws.Range("J131").Value = 747.38776
$ws.Range("L131").Value = 2242.16328
$ws.Range("N131").Value = -12322.16328

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 10000
$ws.Range("I82").Value = 10000
$ws.Range("K82").Value = 10000
$ws.Range("M82").Value = -9617

$ws.Range("H85").Value = 10000
$ws.Range("I85").Value = 10000
$ws.Range("K85").Value = 10000
$ws.Range("M85").Value = -8674

$ws.Range("H113").Value = 2048.3872
$ws.Range("I113").Value = 1568.3334
$ws.Range("J113").Value = 2498.4375
$ws.Range("K113").Value = 1568.3334
$ws.Range("L113").Value = 2498.4375
$ws.Range("M113").Value = 601.6666
$ws.Range("N113").Value = -6838.4375

$ws.Range("H122").Value = 2746.5454
$ws.Range("I122").Value = 2066.6667
$ws.Range("J122").Value = 3562.4
$ws.Range("K122").Value = 6200.000100000001
$ws.Range("L122").Value = 10687.2
$ws.Range("M122").Value = -3750.000100000001
$ws.Range("N122").Value = -15587.2

$ws.Range("H126").Value = 3694.8918
$ws.Range("I126").Value = 2622.6667
$ws.Range("J126").Value = 6589.9
$ws.Range("K126").Value = 7868.000100000001
$ws.Range("L126").Value = 19769.7
$ws.Range("M126").Value = -5398.000100000001
$ws.Range("N126").Value = -24709.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 328.22726
$ws.Range("I16").Value = 332.78946
$ws.Range("K16").Value = 332.78946
$ws.Range("M16").Value = -162.78946

$ws.Range("H21").Value = 1302.5
$ws.Range("I21").Value = 500
$ws.Range("J21").Value = 2105
$ws.Range("K21").Value = 500
$ws.Range("L21").Value = 2105
$ws.Range("M21").Value = -326
$ws.Range("N21").Value = -2453

$ws.Range("H40").Value = 2722.12
$ws.Range("I40").Value = 1745.75
$ws.Range("J40").Value = 3623.3845
$ws.Range("K40").Value = 1745.75
$ws.Range("L40").Value = 3623.3845
$ws.Range("M40").Value = -1609.75
$ws.Range("N40").Value = -3895.3845

$ws.Range("H100").Value = 3551.1
$ws.Range("I100").Value = 1585.5
$ws.Range("J100").Value = 6499.5
$ws.Range("K100").Value = 1585.5
$ws.Range("L100").Value = 6499.5
$ws.Range("M100").Value = -1044.5
$ws.Range("N100").Value = -7581.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 125001160
$ws.Range("I81").Value = 1383.3334
$ws.Range("J81").Value = 500000500
$ws.Range("K81").Value = 2766.6668
$ws.Range("L81").Value = 1000001000
$ws.Range("M81").Value = -1705.6668
$ws.Range("N81").Value = -1000003122

$ws.Range("H84").Value = 125001160
$ws.Range("I84").Value = 1383.3334
$ws.Range("J84").Value = 500000500
$ws.Range("K84").Value = 13833.334
$ws.Range("L84").Value = 5000005000
$ws.Range("M84").Value = -8529.333999999999
$ws.Range("N84").Value = -5000015608

$ws.Range("H113").Value = 5407705.5
$ws.Range("I113").Value = 2875
$ws.Range("J113").Value = 27027028
$ws.Range("K113").Value = 8625
$ws.Range("L113").Value = 81081084
$ws.Range("M113").Value = -6455
$ws.Range("N113").Value = -81085424

$ws.Range("H132").Value = 762
$ws.Range("I132").Value = 762
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2286
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 244
$ws.Range("N132").ClearContents()
